$wb = $excel.ActiveWorkbook

# --- Sheet "09-pressv20": update inputs D27 and D28 (geometry change) ---
$ws1 = $wb.Worksheets.Item("09-pressv20")
$ws1.Range("D27").Value = 17.5
$ws1.Range("D28").Value = 268
[void]$ws1.Range("D51").Select()

# --- Sheet "09-pressv20 (2)": no data changes, just updated selection ---
$ws2 = $wb.Worksheets.Item("09-pressv20 (2)")
[void]$ws2.Range("D51").Select()

# --- Sheet "09-pressv20 (3)": update input D28 (geometry change) ---
$ws3 = $wb.Worksheets.Item("09-pressv20 (3)")
$ws3.Range("D28").Value = 310.2
[void]$ws3.Range("D51").Select()

# --- Sheet "09-pressv20 (4)": update input D28 (geometry change) ---
$ws4 = $wb.Worksheets.Item("09-pressv20 (4)")
$ws4.Range("D28").Value = 275
[void]$ws4.Range("D28").Select()

# --- Sheet "Tabelle4": becomes the active sheet, new selection ---
$ws5 = $wb.Worksheets.Item("Tabelle4")
[void]$ws5.Range("C12").Select()
